# NATMI ligand-receptor (Ngf-Ntrk1) sheet update.
#
# The underlying NATMI run was redone with new TPM values and an updated
# cluster set: the old "Neutrophils" / "Resolving-Mac" clusters are gone and
# a new "ECs" cluster was added, alongside the existing "FAPs" / "MuSCs"
# clusters. That yields a full 3x3 sending x target cluster grid (9 data
# rows instead of the previous 8), each with refreshed NATMI statistics.
#
# Row 1 (headers) is untouched. Rows 2-10 are rewritten in one shot via a
# 2-D array assigned to A2:T10 (faster and more robust than 200 individual
# cell writes), which also naturally grows the sheet's used range /
# dimension from A1:T9 to A1:T10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 9 rows (r=2..10) x 20 columns (A..T)
$data = New-Object 'object[,]' 9,20

# Row 2: sending=ECs, ligand=Ngf, receptor=Ntrk1, target=ECs
$data[0,0] = "ECs"
$data[0,1] = "Ngf"
$data[0,2] = "Ntrk1"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 1
$data[0,6] = 4.5666645
$data[0,7] = 9.133329
$data[0,8] = 0.156583237611307
$data[0,9] = 0.1473423006975575
$data[0,10] = 1
$data[0,11] = 0.5
$data[0,12] = 0.0003185
$data[0,13] = 0.000637
$data[0,14] = 0.001593934541086978
$data[0,15] = 0.001183922444716212
$data[0,16] = 0.00145448264325
$data[0,17] = 0.005817930573
$data[0,18] = 0.0002495834309838918
$data[0,19] = 0.0001744418568519635

# Row 3: sending=ECs, ligand=Ngf, receptor=Ntrk1, target=FAPs
$data[1,0] = "ECs"
$data[1,1] = "Ngf"
$data[1,2] = "Ntrk1"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 1
$data[1,6] = 4.5666645
$data[1,7] = 9.133329
$data[1,8] = 0.156583237611307
$data[1,9] = 0.1473423006975575
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 0.138402
$data[1,13] = 0.415206
$data[1,14] = 0.6926333700330297
$data[1,15] = 0.7716981202210981
$data[1,16] = 0.632035500129
$data[1,17] = 3.792213000774
$data[1,18] = 0.1084547755574022
$data[1,19] = 0.1137037764773569

# Row 4: sending=ECs, ligand=Ngf, receptor=Ntrk1, target=MuSCs
$data[2,0] = "ECs"
$data[2,1] = "Ngf"
$data[2,2] = "Ntrk1"
$data[2,3] = "MuSCs"
$data[2,4] = 2
$data[2,5] = 1
$data[2,6] = 4.5666645
$data[2,7] = 9.133329
$data[2,8] = 0.156583237611307
$data[2,9] = 0.1473423006975575
$data[2,10] = 1
$data[2,11] = 0.5
$data[2,12] = 0.0610995
$data[2,13] = 0.122199
$data[2,14] = 0.3057726954258833
$data[2,15] = 0.2271179573341859
$data[2,16] = 0.27902091761775
$data[2,17] = 1.116083670471
$data[2,18] = 0.04787887862292088
$data[2,19] = 0.03346408236334865

# Row 5: sending=FAPs, ligand=Ngf, receptor=Ntrk1, target=ECs
$data[3,0] = "FAPs"
$data[3,1] = "Ngf"
$data[3,2] = "Ntrk1"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 3.658241666666667
$data[3,7] = 10.974725
$data[3,8] = 0.1254349480088258
$data[3,9] = 0.1770483939670849
$data[3,10] = 1
$data[3,11] = 0.5
$data[3,12] = 0.0003185
$data[3,13] = 0.000637
$data[3,14] = 0.001593934541086978
$data[3,15] = 0.001183922444716212
$data[3,16] = 0.001165149970833333
$data[3,17] = 0.006990899825
$data[3,18] = 0.0001999350962907167
$data[3,19] = 0.0002096115674185902

# Row 6: sending=FAPs, ligand=Ngf, receptor=Ntrk1, target=FAPs
$data[4,0] = "FAPs"
$data[4,1] = "Ngf"
$data[4,2] = "Ntrk1"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 3.658241666666667
$data[4,7] = 10.974725
$data[4,8] = 0.1254349480088258
$data[4,9] = 0.1770483939670849
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.138402
$data[4,13] = 0.415206
$data[4,14] = 0.6926333700330297
$data[4,15] = 0.7716981202210981
$data[4,16] = 0.5063079631499999
$data[4,17] = 4.55677166835
$data[4,18] = 0.08688043075927086
$data[4,19] = 0.1366279128125638

# Row 7: sending=FAPs, ligand=Ngf, receptor=Ntrk1, target=MuSCs
$data[5,0] = "FAPs"
$data[5,1] = "Ngf"
$data[5,2] = "Ntrk1"
$data[5,3] = "MuSCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 3.658241666666667
$data[5,7] = 10.974725
$data[5,8] = 0.1254349480088258
$data[5,9] = 0.1770483939670849
$data[5,10] = 1
$data[5,11] = 0.5
$data[5,12] = 0.0610995
$data[5,13] = 0.122199
$data[5,14] = 0.3057726954258833
$data[5,15] = 0.2271179573341859
$data[5,16] = 0.2235167367125
$data[5,17] = 1.341100420275
$data[5,18] = 0.03835458215326419
$data[5,19] = 0.04021086958710252

# Row 8: sending=MuSCs, ligand=Ngf, receptor=Ntrk1, target=ECs
$data[6,0] = "MuSCs"
$data[6,1] = "Ngf"
$data[6,2] = "Ntrk1"
$data[6,3] = "ECs"
$data[6,4] = 2
$data[6,5] = 1
$data[6,6] = 20.939547
$data[6,7] = 41.87909399999999
$data[6,8] = 0.7179818143798673
$data[6,9] = 0.6756093053353576
$data[6,10] = 1
$data[6,11] = 0.5
$data[6,12] = 0.0003185
$data[6,13] = 0.000637
$data[6,14] = 0.001593934541086978
$data[6,15] = 0.001183922444716212
$data[6,16] = 0.006669245719499999
$data[6,17] = 0.026676982878
$data[6,18] = 0.00114441601381237
$data[6,19] = 0.0007998690204456582

# Row 9: sending=MuSCs, ligand=Ngf, receptor=Ntrk1, target=FAPs
$data[7,0] = "MuSCs"
$data[7,1] = "Ngf"
$data[7,2] = "Ntrk1"
$data[7,3] = "FAPs"
$data[7,4] = 2
$data[7,5] = 1
$data[7,6] = 20.939547
$data[7,7] = 41.87909399999999
$data[7,8] = 0.7179818143798673
$data[7,9] = 0.6756093053353576
$data[7,10] = 2
$data[7,11] = 0.6666666666666666
$data[7,12] = 0.138402
$data[7,13] = 0.415206
$data[7,14] = 0.6926333700330297
$data[7,15] = 0.7716981202210981
$data[7,16] = 2.898075183894
$data[7,17] = 17.388451103364
$data[7,18] = 0.4972981637163567
$data[7,19] = 0.5213664309311773

# Row 10 (new row): sending=MuSCs, ligand=Ngf, receptor=Ntrk1, target=MuSCs
$data[8,0] = "MuSCs"
$data[8,1] = "Ngf"
$data[8,2] = "Ntrk1"
$data[8,3] = "MuSCs"
$data[8,4] = 2
$data[8,5] = 1
$data[8,6] = 20.939547
$data[8,7] = 41.87909399999999
$data[8,8] = 0.7179818143798673
$data[8,9] = 0.6756093053353576
$data[8,10] = 1
$data[8,11] = 0.5
$data[8,12] = 0.0610995
$data[8,13] = 0.122199
$data[8,14] = 0.3057726954258833
$data[8,15] = 0.2271179573341859
$data[8,16] = 1.2793958519265
$data[8,17] = 5.117583407705999
$data[8,18] = 0.2195392346496982
$data[8,19] = 0.1534430053837347

# Write the whole block in one shot (also extends dimension to A1:T10)
$ws.Range("A2:T10").Value2 = $data
